# Commit: "remove initial biomass in the template"
#
# The InitialBiomass column (column C) on the "Components & input parameter"
# sheet is removed. Column H ("RefugeBiomass") used to be computed with
# =ROUND(C3/100,3) (i.e. InitialBiomass/100); once the InitialBiomass column
# is gone that formula can no longer reference it, so the previously
# calculated values are kept as plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Components & input parameter")

# Capture the values produced by the old =ROUND(C3/100,3) formula before we
# blow away the column it depends on.
$refuge3 = $ws.Range("I3").Value2
$refuge4 = $ws.Range("I4").Value2

# Delete column C (InitialBiomass) entirely; D:I shift left to C:H.
$ws.Columns.Item(3).Delete()

# Re-seed the (now plain) RefugeBiomass values in column H (was column I).
$ws.Range("H3").Value = $refuge3
$ws.Range("H4").Value = $refuge4

# --- view-state bookkeeping -------------------------------------------------
# Make "Components & input parameter" the active sheet/tab and select the
# RefugeBiomass cells that were just turned into static values.
$ws.Activate()
$ws.Range("H3:H4").Select()

# The previously active sheet ("Input time-series") is no longer the active
# tab; re-select its original cell so it stops being "tabSelected".
$wsTime = $wb.Worksheets.Item("Input time-series")
$wsTime.Range("A6").Select()

$ws.Activate()
